$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Mean, Var" summary labels (O1 / R1 swapped meaning per diff)
$ws.Range("O1").Value = "0, 0.23"
$ws.Range("R1").Value = "1.1, 0.096"

# Update the underlying p-value inputs; O3:O5 / R3:R5 are formulas that
# recompute automatically from these.
$ws.Range("N3").Value = 0.000000089321010000000006
$ws.Range("Q3").Value = 0.0000033427959999999999

$ws.Range("N4").Value = 0.000001011073
$ws.Range("Q4").Value = 0.000077226799999999997

$ws.Range("N5").Value = 0.000024035519999999999
$ws.Range("Q5").Value = 0.003272409

# Update the view: zoom in, move selection to D6, clear the pinned top-left cell
$ws.Application.ActiveWindow.Zoom = 142
$ws.Range("D6").Select() | Out-Null
